$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# "Pain Control" sheet: capitalise the Yes/No answer columns.
# Column C (row 3-22) holds "no" -> "No" for every question.
# Column B (row 8-22) holds "yes" -> "Yes" (rows 3-7 hold age-range
# text in column B and are left untouched).
# ------------------------------------------------------------------
$painControl = $wb.Worksheets.Item("Pain Control")

for ($r = 3; $r -le 22; $r++) {
    if ($painControl.Cells.Item($r, 3).Value() -eq "no") {
        $painControl.Cells.Item($r, 3).Value = "No"
    }
    if ($painControl.Cells.Item($r, 2).Value() -eq "yes") {
        $painControl.Cells.Item($r, 2).Value = "Yes"
    }
}

# New blank, but formatted, rows appended below the table (23-27).
$painControl.Range("B23:B27").Font.Name = "Arial"
$painControl.Range("B23:B27").Font.Size = 10
$painControl.Range("B23:B27").Font.Bold = $false
$painControl.Range("B23:B27").Font.Color = 0

$painControl.Rows.Item(24).RowHeight = 15.75
$painControl.Rows.Item(25).RowHeight = 15.75
$painControl.Rows.Item(26).RowHeight = 15.75
$painControl.Rows.Item(27).RowHeight = 15.75

# Cursor left at B34 on this sheet.
[void]$painControl.Range("B34").Select()

# ------------------------------------------------------------------
# "Allergies" sheet: selection moved from B6 to B8:B17.
# ------------------------------------------------------------------
$allergies = $wb.Worksheets.Item("Allergies")
[void]$allergies.Range("B8:B17").Select()

# ------------------------------------------------------------------
# GERD becomes the active/visible tab (was Allergies before).
# ------------------------------------------------------------------
$gerd = $wb.Worksheets.Item("GERD")
$gerd.Activate()
